$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active/selected so we don't disturb it
$aboutSheet = $wb.Worksheets.Item("About")

# The data sheet is "BDSBaPCF"
$ws = $wb.Worksheets.Item("BDSBaPCF")

# Set petroleum's "Do Suppliers Bid at Peak Capacity Factors" flag to TRUE (1)
$ws.Range("B9").Value = 1

# Update the (inactive) selection on BDSBaPCF to match the saved state in the source file
$ws.Activate()
$ws.Range("B5").Select()

# Restore the originally active sheet ("About") as the selected tab
$aboutSheet.Activate()

$wb.Save()
